# Purchase details ledger update: 17/10/2017 MAMATHA CHICK IN
#
# 1) Merge the two runs that make up "SAT Oct 14 14:26:46 PDT 2017" into a
#    single run (cosmetic normalisation left over from a prior edit).
# 2) Append a brand-new ledger entry (MON Oct 16 ... Amount balance) right
#    after the existing last entry (the one ending "Amount balance - 225912.0").

$d = $word.ActiveDocument

# --- 1) Merge "SAT Oct 14" / " 14:26:46 PDT 2017" runs -----------------
$null = $d.Content.Find.Execute(
    "SAT Oct 14 14:26:46 PDT 2017", $true, $false, $false, $false, $false,
    $true, 1, $false, "SAT Oct 14 14:26:46 PDT 2017", 2)

# --- 2) Insert the new "MON Oct 16" ledger entry ------------------------

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-RunProps([bool]$Bold) {
    $bTag = ""
    if ($Bold) { $bTag = "<w:b/>" }
    return "<w:rPr><w:rFonts w:ascii='Courier New' w:hAnsi='Courier New' w:cs='Courier New'/>$bTag</w:rPr>"
}

function New-Run([string]$Text, [bool]$Tab, [bool]$Bold) {
    $xml = "<w:r>" + (Get-RunProps $Bold)
    if ($Tab) { $xml += "<w:tab/>" }
    if ($null -ne $Text) {
        if ($Text.StartsWith(" ") -or $Text.EndsWith(" ")) {
            $xml += "<w:t xml:space='preserve'>$Text</w:t>"
        } else {
            $xml += "<w:t>$Text</w:t>"
        }
    }
    $xml += "</w:r>"
    return $xml
}

function New-Para([string]$RunsXml, [bool]$Bold) {
    $pPr = "<w:pPr><w:pStyle w:val='PlainText'/>" + (Get-RunProps $Bold) + "</w:pPr>"
    return "<w:p $wNs>$pPr$RunsXml</w:p>"
}

# A "label [tabs] - value" row, e.g. "Person Name` t` t` t` t- NG"
function New-FieldPara([string]$Label, [int]$TabCount, [string]$Value, [bool]$Bold) {
    $runs = New-Run $Label $false $Bold
    for ($i = 0; $i -lt $TabCount; $i++) {
        $runs += New-Run $null $true $Bold
    }
    $runs += New-Run $Value $true $Bold
    return New-Para $runs $Bold
}

$xmlBlock = ""

# Blank bold separator line (matches the blank line already preceding each entry)
$xmlBlock += New-Para "" $true

# Timestamp line: two runs, "MON Oct 16" + " 14:14:45 PDT 2017"
$timestampRuns = (New-Run "MON Oct 16" $false $false) + (New-Run " 14:14:45 PDT 2017" $false $false)
$xmlBlock += New-Para $timestampRuns $false

$xmlBlock += New-FieldPara "Person Name" 3 "- NG" $false
$xmlBlock += New-Para (New-Run "---------------------------------------------------------------" $false $false) $false
$xmlBlock += New-FieldPara "Item Name" 3 "- CARROT" $false
$xmlBlock += New-FieldPara "Number of Pockets" 2 "- 5" $false
$xmlBlock += New-FieldPara "Number of KGs" 2 "- 440" $false
$xmlBlock += New-FieldPara "Rate" 4 "- 25" $false

# Transport & Miscellaneous: label immediately followed by a single tab+value run
$transportRuns = (New-Run "Transport &amp; Miscellaneous" $false $false) + (New-Run "- 75" $true $false)
$xmlBlock += New-Para $transportRuns $false

$xmlBlock += New-FieldPara "Total Price" 3 "- 11075.0" $false
$xmlBlock += New-FieldPara "Amount balance" 2 "- 236987.0" $true

# Two trailing blank bold paragraphs
$xmlBlock += New-Para "" $true
$xmlBlock += New-Para "" $true

# Locate the last entry's "Amount balance ... - 225912.0" paragraph and insert
# the new block immediately after it.
$searchRange = $d.Content
$null = $searchRange.Find.Execute("- 225912.0")
$targetPara = $searchRange.Paragraphs(1)
$insertPoint = $d.Range($targetPara.Range.End, $targetPara.Range.End)
$insertPoint.InsertXML($xmlBlock)
